# Apply updated Joule (B column) measurement values for the Nexus 5X Cpu sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$newValues = @{
    2 = 122.01944399999999
    3 = 129.537972
    4 = 123.585804
    5 = 124.838892
    6 = 124.055712
    7 = 126.71852399999899
    8 = 127.65834
    9 = 123.115895999999
    10 = 124.21234800000001
    11 = 124.055712
    12 = 122.95926
    13 = 119.82653999999999
    14 = 124.055712
    15 = 123.89907599999999
    16 = 124.682256
    17 = 128.75479200000001
    18 = 123.115895999999
    19 = 123.429167999999
    20 = 124.368984
    21 = 124.52562
    22 = 123.115895999999
    23 = 123.429167999999
    24 = 124.055712
    25 = 127.188431999999
    26 = 122.489352
    27 = 123.89907599999999
    28 = 122.489352
    29 = 123.585804
    30 = 119.043359999999
    31 = 123.429167999999
    32 = 127.814976
    33 = 134.70695999999899
    34 = 128.59815599999999
    35 = 129.537972
    36 = 127.97161199999999
    37 = 129.537972
    38 = 128.12824800000001
    39 = 128.59815599999999
    40 = 128.44152
    41 = 128.59815599999999
    42 = 129.537972
    43 = 124.838892
    44 = 136.89986399999901
    45 = 128.59815599999999
    46 = 129.22469999999899
    47 = 128.12824800000001
    48 = 128.28488399999901
    49 = 127.814976
    50 = 129.22469999999899
    51 = 128.75479200000001
    52 = 132.98396399999999
    53 = 128.44152
    54 = 129.537972
    55 = 131.104332
    56 = 130.32115200000001
    57 = 134.23705200000001
    58 = 129.06806399999999
    59 = 133.297236
    60 = 124.52562
    61 = 129.22469999999899
    62 = 144.26175599999999
    63 = 145.82811599999999
    64 = 144.88829999999999
    65 = 145.35820799999999
    66 = 147.394476
    67 = 144.26175599999999
    68 = 144.105119999999
    69 = 144.26175599999999
    70 = 143.635212
    71 = 146.92456799999999
    72 = 145.67148
    73 = 142.69539599999999
    74 = 145.201572
    75 = 144.575028
    76 = 145.35820799999999
    77 = 144.73166399999999
    78 = 146.14138800000001
    79 = 150.213923999999
    80 = 144.26175599999999
    81 = 143.635212
    82 = 144.26175599999999
    83 = 143.79184799999999
    84 = 143.635212
    85 = 145.201572
    86 = 145.04493600000001
    87 = 145.04493600000001
    88 = 144.105119999999
    89 = 143.478576
    90 = 144.26175599999999
    91 = 144.41839200000001
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item([int]$row, 2).Value = $newValues[$row]
}

# Restore the view state captured in the saved workbook (active cell / scroll position).
$ws.Activate()
$ws.Range("I76").Select()
$excel.ActiveWindow.ScrollRow = 67
